# Refresh the cryptocurrency price/volume snapshot (Price column D, Volume(1h) column E)
# to the new scraped values. Numeric-looking Price strings are briefly forced to text
# (NumberFormat "@") so Excel doesn't silently coerce them into numbers / lose trailing
# zeros, then ClearFormats() drops the temporary formatting so no stray style lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.510.63'
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = '1.955.05'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.25'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("E6").Value = '  +2.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.08'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +5.06%  '

$ws.Range("E9").Value = '  +3.63%  '

$ws.Range("E10").Value = '  -2.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.845'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.06'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.81%  '

$ws.Range("D14").Value = '2.242.95'
$ws.Range("E14").Value = '  +0.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.55'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.25'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.83%  '

$ws.Range("D17").Value = '1.955.06'
$ws.Range("E17").Value = '  +0.79%  '

$ws.Range("D18").Value = '36.424.43'
$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.09'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.20%  '

$ws.Range("D20").Value = '0.0₃0853'
$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.24'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.06'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.09%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.34%  '

$ws.Range("E26").Value = '  +6.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.09'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.67'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.20'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.30'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +19.82%  '

$ws.Range("E31").Value = '  +2.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.78'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.23%  '

$ws.Range("E33").Value = '  -0.85%  '

$ws.Range("E34").Value = '  +6.62%  '

$ws.Range("E35").Value = '  +0.17%  '

$ws.Range("E36").Value = '  +2.75%  '

$ws.Range("E37").Value = '  +3.23%  '

$ws.Range("E38").Value = '  -0.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.41'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -11.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0964'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.08%  '

$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("E42").Value = '  +1.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0209'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.81'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.43%  '

$ws.Range("D45").Value = '1.359.52'
$ws.Range("E45").Value = '  +1.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.48'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.66%  '

$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.21'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.23%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.87'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +6.63%  '

$ws.Range("D51").Value = '2.138.37'
$ws.Range("E51").Value = '  +0.94%  '
